# "Generate Report for Handoff"
#
# For the four files that are "Ready for handoff" (rows 4-7 on both the
# zh-cn and de-de localization-status sheets), mark them as handed off by
# bumping their Priority to "ht" and stamping a fresh "Latest Handoff
# Datetime".

$wb = $excel.ActiveWorkbook

$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

$zhHandoffTime = "2016-08-27 06:29:13"
$deHandoffTime = "2016-08-27 06:29:18"

for ($row = 4; $row -le 7; $row++) {
    # Column E = Priority, Column H = Latest Handoff Datetime
    $zhSheet.Range("E" + $row).Value = "ht"
    $zhSheet.Range("H" + $row).Value = $zhHandoffTime

    $deSheet.Range("E" + $row).Value = "ht"
    $deSheet.Range("H" + $row).Value = $deHandoffTime
}
